$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.447.84"
$ws.Range("E2").Value = "  +0.57%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.682.94"
$ws.Range("E3").Value = "  +0.19%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "686.97"
$ws.Range("E5").Value = "  +1.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.88"
$ws.Range("E6").Value = "  -0.91%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.26%  "

# Row 8
$ws.Range("E8").Value = "  -0.20%  "

# Row 9
$ws.Range("E9").Value = "  -1.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.08"
$ws.Range("E10").Value = "  -2.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("E11").Value = "  -3.28%  "

# Row 12
$ws.Range("E12").Value = "  -0.78%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.304.34"
$ws.Range("E13").Value = "  +0.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.26"
$ws.Range("E14").Value = "  -2.84%  "

# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.438.47"
$ws.Range("E15").Value = "  +0.43%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.670.53"
$ws.Range("E16").Value = "  -0.14%  "

# Row 17
$ws.Range("E17").Value = "  +2.03%  "

# Row 18
$ws.Range("E18").Value = "  -2.56%  "

# Row 19
$ws.Range("E19").Value = "  -3.37%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.66"
$ws.Range("E20").Value = "  -2.53%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.95"
$ws.Range("E21").Value = "  +2.09%  "

# Row 22
$ws.Range("E22").Value = "  -1.55%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.60"
$ws.Range("E23").Value = "  +0.31%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.828.80"
$ws.Range("E24").Value = "  +0.18%  "

# Row 26
$ws.Range("E26").Value = "  -2.52%  "

# Row 27
$ws.Range("E27").Value = "  -4.73%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.17"
$ws.Range("E28").Value = "  -3.12%  "

# Row 29
$ws.Range("E29").Value = "  -0.65%  "

# Row 30
$ws.Range("E30").Value = "  -4.67%  "

# Row 31
$ws.Range("E31").Value = "  -4.70%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.57"
$ws.Range("E32").Value = "  -1.89%  "

# Row 33
$ws.Range("E33").Value = "  -0.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.80"
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.655.72"
$ws.Range("E35").Value = "  +0.37%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.160"
$ws.Range("E36").Value = "  -2.15%  "

# Row 37
$ws.Range("E37").Value = "  -3.44%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.12"
$ws.Range("E38").Value = "  +1.53%  "

# Row 39
$ws.Range("E39").Value = "  +0.01%  "

# Row 40
$ws.Range("E40").Value = "  +2.50%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0897"
$ws.Range("E41").Value = "  -4.04%  "

# Row 42
$ws.Range("E42").Value = "  -0.03%  "

# Row 43
$ws.Range("E43").Value = "  -1.40%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "165.58"
$ws.Range("E44").Value = "  +5.28%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.50"
$ws.Range("E45").Value = "  -0.98%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000282"
$ws.Range("E46").Value = "  +1.99%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.72"
$ws.Range("E47").Value = "  -2.42%  "

# Row 48
$ws.Range("E48").Value = "  +5.90%  "

# Row 49
$ws.Range("E49").Value = "  +0.28%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.71"
$ws.Range("E50").Value = "  -1.42%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.78"
$ws.Range("E51").Value = "  -2.81%  "
